{"js": "// Remove the trailing \"Ver no Jupiter ...\" / copyright boilerplate block\n// (and the blank paragraph that preceded it) that the site generator used\n// to append to each course page, while leaving the rest of the document\n// (including the final blank paragraph + page-break paragraph) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\n\n// Locate the contiguous run: [blank paragraph] \"Ver no Jupiter...\" \"\u00a9 2020...\"\nlet blockStart = -1;\nlet blockEnd = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (targetTexts.indexOf(items[i].text) !== -1) {\n    if (blockStart === -1) {\n      blockStart = i;\n    }\n    blockEnd = i;\n  }\n}\n\nif (blockStart !== -1) {\n  // Include the blank paragraph immediately preceding the block, if present.\n  if (blockStart > 0 && items[blockStart - 1].text === \"\") {\n    blockStart = blockStart - 1;\n  }\n\n  // Delete from the bottom up so earlier indices stay valid.\n  for (let i = blockEnd; i >= blockStart; i--) {\n    items[i].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / copyright boilerplate block\n# (and the blank paragraph that preceded it) that the site generator used\n# to append to each course page, while leaving the rest of the document\n# (including the final blank paragraph + page-break paragraph) untouched.\n\n$d = $word.ActiveDocument\n\n$targetTexts = @(\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n# Find the paragraph indices (1-based) of the boilerplate block.\n$blockStart = -1\n$blockEnd = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($targetTexts -contains $t) {\n        if ($blockStart -eq -1) {\n            $blockStart = $i\n        }\n        $blockEnd = $i\n    }\n    $i++\n}\n\nif ($blockStart -ne -1) {\n    # Also remove the blank paragraph immediately preceding the block, if any.\n    if ($blockStart -gt 1) {\n        $prev = $d.Paragraphs.Item($blockStart - 1)\n        $prevText = $prev.Range.Text.TrimEnd([char]13, [char]7)\n        if ($prevText -eq \"\") {\n            $blockStart = $blockStart - 1\n        }\n    }\n\n    # Delete bottom-up so earlier paragraph indices stay valid.\n    for ($idx = $blockEnd; $idx -ge $blockStart; $idx--) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
